$p = $ppt.ActivePresentation

# --- Slide 1: rotate the team member name list up by one position ---
# (VILMARD Alexis moves from first to last; a new trailing blank line is added)
$s1 = $p.Slides.Item(1)
$namesShape = $s1.Shapes.Item(5)
$namesTr = $namesShape.TextFrame.TextRange

$newNames = @("BERNHARD William", "GUILLOU Aurélien", "HEIDET Lucas", "TROGNOT Mathias", "VILMARD Alexis")

for ($i = 1; $i -le 5; $i++) {
    $para = $namesTr.Paragraphs($i)
    # Go through a throwaway value first so the run is fully replaced instead
    # of being diffed/split against the previous name (they share letters).
    $para.Text = "."
    $para.Text = $newNames[$i - 1]
}

# Add the new empty trailing paragraph (keeps the same paragraph formatting).
[void]$namesTr.InsertAfter("`r")

# --- Slide 2: fix the gender agreement "choisie" -> "choisi" ("jeu" is masculine) ---
$s2 = $p.Slides.Item(2)
$introShape = $s2.Shapes.Item(2)
$introTr = $introShape.TextFrame.TextRange

$introPara1 = $introTr.Paragraphs(1)

# "Présentation du jeu choisie" -> "Présentation " + "du jeu " + "choisi"
$middle = $introTr.Characters(14, 7)
$middle.Text = "du jeu "

$suffix = $introTr.Characters(21, 7)
$suffix.Text = "choisi"
